$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: UserID 1001 / A / 30000 / 2025-02 / <hash>  ->  maythawee / A / 30000 / 2025-02 / Maymys@393833
$ws.Range("A2").Value = "maythawee"

# --- Header row: HashedPassword -> Password --------------------------------
$ws.Range("E1").Value = "Password"

$ws.Range("B2").Value = "A"
$ws.Range("C2").Value = 30000
$ws.Range("D2").Value = "2025-02"
$ws.Range("E2").Value = "Maymys@393833"

# --- Row 3: UserID 1002 / B / 35000 / 2025-02 / <hash> -> admin / B / 35000 / 2025-02 / Admin@393833
$ws.Range("A3").Value = "admin"
$ws.Range("B3").Value = "B"
$ws.Range("C3").Value = 35000
$ws.Range("D3").Value = "2025-02"
$ws.Range("E3").Value = "Admin@393833"

# --- Row 4 (UserID 1003 / Admin / 0 / - / <hash>) is gone entirely --------
$ws.Rows("4:4").Delete()

# --- Give the data block a thin border (A2:D3) -----------------------------
$ws.Range("A2:D3").Borders.LineStyle = 1

# --- Border the Password column too, then turn the two values into real
#     hyperlinks (Excel auto-links "name@host"-shaped text) ----------------
$ws.Range("E2:E3").Borders.LineStyle = 1
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Maymys@393833")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:Admin@393833")

# --- Resize the columns: A-D auto-fit to their (now short) content, E gets
#     a fixed, human width instead of the old 164-char bestFit -------------
$ws.Columns("A:A").AutoFit()
$ws.Columns("B:B").AutoFit()
$ws.Columns("C:C").AutoFit()
$ws.Columns("D:D").AutoFit()
$ws.Columns("E:E").ColumnWidth = 15.7

# --- Selection moved off the old E22 default -------------------------------
$ws.Range("E10").Select() | Out-Null
